$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.095.08"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.285.72"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.584"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.81%  "
$ws.Range("D9").Value = "2.283.25"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0996"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "2.689.88"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "58.003.20"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "2.277.03"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("E28").Value = "  -2.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("D31").Value = "0.0₃0722"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.378"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "140.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "287.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0952"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.552"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0210"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.62%  "
